$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a value into a cell while forcing it to be stored as TEXT,
# even when the value "looks like" a number or a date (e.g. "14316317" or
# "2025-08-11"), which a plain Range.Value assignment would otherwise let
# Excel auto-convert into a Double / date serial.
#
# We stage the text in a scratch cell as a formula that evaluates to a text
# string ( ="14316317" ), which always yields a text result regardless of
# how the literal looks; then Copy + PasteSpecial (values only) the computed
# text into the destination. The destination cell's NumberFormat/style is
# never touched, so it keeps the workbook's default (General) style while
# still holding a genuine text value - matching how the source rows were
# authored.
function Set-TextValue {
    param($targetRange, [string]$text)
    $scratch = $ws.Range("Z1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $targetRange.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $excel.CutCopyMode = 0
    $scratch.Clear() | Out-Null
}

# New "tracker_resultados" rows appended by the automatic 3h sync.
$newRows = @(
    @{ Row=163; EventId="14316317"; Fecha="2025-08-11"; JugA="Frances Tiafoe";         JugB="Ugo Humbert";     Pron="Gana Frances Tiafoe";         Cuota=1.67 },
    @{ Row=164; EventId="14316318"; Fecha="2025-08-11"; JugA="Joao Fonseca";           JugB="Terence Atmane";  Pron="Gana Terence Atmane";         Cuota=4    },
    @{ Row=165; EventId="14316410"; Fecha="2025-08-11"; JugA="Jessica Bouzas Maneiro"; JugB="Taylor Townsend"; Pron="Gana Jessica Bouzas Maneiro"; Cuota=2.1  },
    @{ Row=166; EventId="14316446"; Fecha="2025-08-11"; JugA="Aryna Sabalenka";        JugB="Emma Raducanu";   Pron="Gana Aryna Sabalenka";        Cuota=1.33 },
    @{ Row=167; EventId="14393238"; Fecha="2025-08-11"; JugA="Stefano Travaglia";      JugB="Tiago Pereira";   Pron="Gana Tiago Pereira";          Cuota=1.83 },
    @{ Row=168; EventId="14316317"; Fecha="2025-08-11"; JugA="Frances Tiafoe";         JugB="Ugo Humbert";     Pron="Gana Frances Tiafoe";         Cuota=1.67 },
    @{ Row=169; EventId="14316318"; Fecha="2025-08-11"; JugA="Joao Fonseca";           JugB="Terence Atmane";  Pron="Gana Terence Atmane";         Cuota=4    }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    Set-TextValue $ws.Range("A$rowNum") $r.EventId   # event_id, kept as text (matches source feed)
    Set-TextValue $ws.Range("B$rowNum") $r.Fecha      # fecha, kept as text (not an Excel date)
    $ws.Range("C$rowNum").Value = $r.JugA             # jugador_A
    $ws.Range("D$rowNum").Value = $r.JugB             # jugador_B
    $ws.Range("E$rowNum").Value = $r.Pron             # pronostico
    $ws.Range("F$rowNum").Value = $r.Cuota            # cuota
    # resultado / profit (G, H) are left blank - the match hasn't been
    # settled yet, same as the other still-pending rows in the tracker.
}

Write-Output ("Added rows " + $newRows[0].Row + "-" + $newRows[-1].Row)
